# Reorder the comma-separated "Recorded By" names/emails in column G.
# For each non-empty, multi-value cell in column G (except the literal
# value "System, admin@admin.com", which the source diff leaves untouched),
# reverse the order of the comma-separated entries.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$used = $ws.UsedRange
$lastRow = $used.Row + $used.Rows.Count - 1

for ($r = 1; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)
    $val = $cell.Value2

    if ($val -eq $null) {
        continue
    }

    if ($val -eq "") {
        continue
    }

    if ($val.IndexOf(",") -lt 0) {
        continue
    }

    if ($val -eq "System, admin@admin.com") {
        continue
    }

    $parts = $val.Split(",")
    $trimmed = @()
    foreach ($p in $parts) {
        $trimmed += $p.Trim()
    }

    $n = $trimmed.Length
    $rev = @()
    for ($i = $n - 1; $i -ge 0; $i--) {
        $rev += $trimmed[$i]
    }

    $newVal = [string]::Join(", ", $rev)
    $cell.Value = $newVal
}
